$d = $word.ActiveDocument
$found = $d.Content.Find.Execute("Dheeraj Chand", $true, $false, $false, $false, $false, $true, 1, $false, "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX", 2)
Write-Output "found=$found"
for ($i = 1; $i -le 3; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output "$i : $($p.Range.Text)"
}
